$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Icam1", "Itgam", "ECs", 3, 1, 144.843106, 434.529318, 0.6517202749316883, 0.6517202749316884, 1, 0.3333333333333333, 66.38494866666666, 199.154846, 0.4171200956172241, 0.4171200956172242, 9615.402156530557, 86538.61940877502, 0.2718456233951894, 0.2718456233951895),
    @("ECs", "Icam1", "Itgam", "M2", 3, 1, 144.843106, 434.529318, 0.6517202749316883, 0.6517202749316884, 3, 1, 92.76573566666667, 278.297207, 0.5828799043827758, 0.5828799043827758, 13436.47728433498, 120928.2955590148, 0.3798746515364988, 0.3798746515364989),
    @("FAPs", "Icam1", "Itgam", "ECs", 3, 1, 33.793597, 101.380791, 0.1520539909422453, 0.1520539909422453, 1, 0.3333333333333333, 66.38494866666666, 199.154846, 0.4171200956172241, 0.4171200956172242, 2243.38620210702, 20190.47581896319, 0.06342477524080989, 0.0634247752408099),
    @("FAPs", "Icam1", "Itgam", "M2", 3, 1, 33.793597, 101.380791, 0.1520539909422453, 0.1520539909422453, 3, 1, 92.76573566666667, 278.297207, 0.5828799043827758, 0.5828799043827758, 3134.88788652786, 28213.99097875074, 0.08862921570143541, 0.08862921570143541),
    @("M2", "Icam1", "Itgam", "ECs", 3, 1, 41.33760066666667, 124.012802, 0.1859981687460937, 0.1859981687460937, 1, 0.3333333333333333, 66.38494866666666, 199.154846, 0.4171200956172241, 0.4171200956172242, 2744.194498259832, 24697.75048433849, 0.0775835739319992, 0.07758357393199922),
    @("M2", "Icam1", "Itgam", "M2", 3, 1, 41.33760066666667, 124.012802, 0.1859981687460937, 0.1859981687460937, 3, 1, 92.76573566666667, 278.297207, 0.5828799043827758, 0.5828799043827758, 3834.712936538224, 34512.41642884402, 0.1084145948140945, 0.1084145948140945),
    @("sCs", "Icam1", "Itgam", "ECs", 3, 1, 2.273049333333333, 6.819148, 0.01022756537997252, 0.01022756537997253, 1, 0.3333333333333333, 66.38494866666666, 199.154846, 0.4171200956172241, 0.4171200956172242, 150.8962633101342, 1358.066369791208, 0.004266123049225551, 0.004266123049225552),
    @("sCs", "Icam1", "Itgam", "M2", 3, 1, 2.273049333333333, 6.819148, 0.01022756537997252, 0.01022756537997253, 3, 1, 92.76573566666667, 278.297207, 0.5828799043827758, 0.5828799043827758, 210.8610936132929, 1897.749842519636, 0.005961442330746973, 0.005961442330746975),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}
